# Add 3 new price records (rows) for Papa (potato) at the top of the
# "Macroferia Regional de Talca" weekly block, shifting the existing
# records starting at row 346 down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows at row 346 (existing rows 346.. shift to 349..)
$ws.Rows.Item(346).Insert()
$ws.Rows.Item(346).Insert()
$ws.Rows.Item(346).Insert()

# ---- New row 346: Asterix, 1a (cosecha) ----
$ws.Cells.Item(346,1).Value  = 5
$ws.Cells.Item(346,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(346,3).Value  = "Maule"
$ws.Cells.Item(346,4).Value  = 44642
$ws.Cells.Item(346,5).Value  = 7
$ws.Cells.Item(346,6).Value  = 100114001
$ws.Cells.Item(346,7).Value  = "Papa"
$ws.Cells.Item(346,8).Value  = "Asterix"
$ws.Cells.Item(346,9).Value  = "1a (cosecha)"
$ws.Cells.Item(346,10).Value = 1200
$ws.Cells.Item(346,11).Value = 7000
$ws.Cells.Item(346,12).Value = 7000
$ws.Cells.Item(346,13).Value = 7000
$ws.Cells.Item(346,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(346,15).Value = "Región de Los Lagos"
$ws.Cells.Item(346,16).Value = 280
$ws.Cells.Item(346,17).Value = 25
$ws.Cells.Item(346,18).Value = "Hortaliza"

# ---- New row 347: Patagonia, 1a (cosecha) ----
$ws.Cells.Item(347,1).Value  = 5
$ws.Cells.Item(347,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(347,3).Value  = "Maule"
$ws.Cells.Item(347,4).Value  = 44642
$ws.Cells.Item(347,5).Value  = 7
$ws.Cells.Item(347,6).Value  = 100114001
$ws.Cells.Item(347,7).Value  = "Papa"
$ws.Cells.Item(347,8).Value  = "Patagonia"
$ws.Cells.Item(347,9).Value  = "1a (cosecha)"
$ws.Cells.Item(347,10).Value = 1200
$ws.Cells.Item(347,11).Value = 6500
$ws.Cells.Item(347,12).Value = 6500
$ws.Cells.Item(347,13).Value = 6500
$ws.Cells.Item(347,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(347,15).Value = "Región de Los Lagos"
$ws.Cells.Item(347,16).Value = 260
$ws.Cells.Item(347,17).Value = 25
$ws.Cells.Item(347,18).Value = "Hortaliza"

# ---- New row 348: Yagana, 1a (cosecha) ----
$ws.Cells.Item(348,1).Value  = 5
$ws.Cells.Item(348,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(348,3).Value  = "Maule"
$ws.Cells.Item(348,4).Value  = 44642
$ws.Cells.Item(348,5).Value  = 7
$ws.Cells.Item(348,6).Value  = 100114001
$ws.Cells.Item(348,7).Value  = "Papa"
$ws.Cells.Item(348,8).Value  = "Yagana"
$ws.Cells.Item(348,9).Value  = "1a (cosecha)"
$ws.Cells.Item(348,10).Value = 800
$ws.Cells.Item(348,11).Value = 8000
$ws.Cells.Item(348,12).Value = 8000
$ws.Cells.Item(348,13).Value = 8000
$ws.Cells.Item(348,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(348,15).Value = "Región de Los Lagos"
$ws.Cells.Item(348,16).Value = 320
$ws.Cells.Item(348,17).Value = 25
$ws.Cells.Item(348,18).Value = "Hortaliza"
